# chore: update Sheets via scheduled runner
#
# Refresh cached Market Board price / leve-profit figures (columns
# H:N = currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across the crafting-class sheets. Values only —
# no structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 643.375
$ws.Range("I53").Value = 415
$ws.Range("J53").Value = 780.4
$ws.Range("K53").Value = 415
$ws.Range("L53").Value = 780.4
$ws.Range("M53").Value = 222
$ws.Range("N53").Value = -2054.4

$ws.Range("H55").Value = 172.72728
$ws.Range("I55").Value = 128.57143
$ws.Range("J55").Value = 250
$ws.Range("K55").Value = 128.57143
$ws.Range("L55").Value = 250
$ws.Range("M55").Value = 85.42857000000001
$ws.Range("N55").Value = -678

$ws.Range("H80").Value = 1579.8
$ws.Range("I80").Value = 1133
$ws.Range("K80").Value = 3399
$ws.Range("M80").Value = -2401

$ws.Range("H83").Value = 1579.8
$ws.Range("I83").Value = 1133
$ws.Range("K83").Value = 10197
$ws.Range("M83").Value = -5205

$ws.Range("H111").Value = 808
$ws.Range("I111").Value = 717
$ws.Range("J111").Value = 990
$ws.Range("K111").Value = 2151
$ws.Range("L111").Value = 2970
$ws.Range("M111").Value = 916
$ws.Range("N111").Value = -9104

$ws.Range("H116").Value = 376840.75
$ws.Range("I116").Value = 910909.0600000001
$ws.Range("J116").Value = 9668.75
$ws.Range("K116").Value = 910909.0600000001
$ws.Range("L116").Value = 9668.75
$ws.Range("M116").Value = -907467.0600000001
$ws.Range("N116").Value = -16552.75

$ws.Range("H118").Value = 510.44446
$ws.Range("I118").Value = 273.42856
$ws.Range("K118").Value = 820.28568
$ws.Range("M118").Value = 836.71432

$ws.Range("H138").Value = 3121.7124
$ws.Range("I138").Value = 1798.1428
$ws.Range("J138").Value = 3435.7795
$ws.Range("K138").Value = 5394.428400000001
$ws.Range("L138").Value = 10307.3385
$ws.Range("M138").Value = -254.4284000000007
$ws.Range("N138").Value = -20587.3385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 887.0857
$ws.Range("I2").Value = 878.96155
$ws.Range("J2").Value = 910.55554
$ws.Range("K2").Value = 878.96155
$ws.Range("L2").Value = 910.55554
$ws.Range("M2").Value = -765.96155
$ws.Range("N2").Value = -1136.55554

$ws.Range("H32").Value = 7266.681
$ws.Range("I32").Value = 5556.0464
$ws.Range("K32").Value = 5556.0464
$ws.Range("M32").Value = -5269.0464

$ws.Range("H63").Value = 6299249.5
$ws.Range("I63").Value = 9896142
$ws.Range("K63").Value = 9896142
$ws.Range("M63").Value = -9895456

$ws.Range("H66").Value = 6299249.5
$ws.Range("I66").Value = 9896142
$ws.Range("K66").Value = 49480710
$ws.Range("M66").Value = -49477278

$ws.Range("H116").Value = 887.0857
$ws.Range("I116").Value = 878.96155
$ws.Range("J116").Value = 910.55554
$ws.Range("K116").Value = 878.96155
$ws.Range("L116").Value = 910.55554
$ws.Range("M116").Value = 1415.03845
$ws.Range("N116").Value = -5498.55554

$ws.Range("H137").Value = 39785.715
$ws.Range("J137").Value = 39785.715
$ws.Range("L137").Value = 39785.715
$ws.Range("N137").Value = -49985.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 887.0857
$ws.Range("I3").Value = 878.96155
$ws.Range("J3").Value = 910.55554
$ws.Range("K3").Value = 878.96155
$ws.Range("L3").Value = 910.55554
$ws.Range("M3").Value = -764.96155
$ws.Range("N3").Value = -1138.55554

$ws.Range("H59").Value = 118830
$ws.Range("J59").Value = 118830
$ws.Range("L59").Value = 118830
$ws.Range("N59").Value = -120524

$ws.Range("H94").Value = 924.8148
$ws.Range("I94").Value = 1000.6429
$ws.Range("J94").Value = 843.1539
$ws.Range("K94").Value = 1000.6429
$ws.Range("L94").Value = 843.1539
$ws.Range("M94").Value = -549.6429000000001
$ws.Range("N94").Value = -1745.1539

$ws.Range("H99").Value = 3829
$ws.Range("I99").Value = 1148
$ws.Range("J99").Value = 4786.5
$ws.Range("K99").Value = 1148
$ws.Range("L99").Value = 4786.5
$ws.Range("M99").Value = 350
$ws.Range("N99").Value = -7782.5

$ws.Range("H137").Value = 36225
$ws.Range("J137").Value = 39966.668
$ws.Range("L137").Value = 39966.668
$ws.Range("N137").Value = -50166.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12198737
$ws.Range("I31").Value = 2214.96
$ws.Range("J31").Value = 31255804
$ws.Range("K31").Value = 2214.96
$ws.Range("L31").Value = 31255804
$ws.Range("M31").Value = -1919.96
$ws.Range("N31").Value = -31256394

$ws.Range("H34").Value = 12198737
$ws.Range("I34").Value = 2214.96
$ws.Range("J34").Value = 31255804
$ws.Range("K34").Value = 2214.96
$ws.Range("L34").Value = 31255804
$ws.Range("M34").Value = -2012.96
$ws.Range("N34").Value = -31256208

$ws.Range("H99").Value = 18186702
$ws.Range("I99").Value = 28573818
$ws.Range("J99").Value = 9250
$ws.Range("K99").Value = 28573818
$ws.Range("L99").Value = 9250
$ws.Range("M99").Value = -28572320
$ws.Range("N99").Value = -12246

$ws.Range("H126").Value = 18186702
$ws.Range("I126").Value = 28573818
$ws.Range("J126").Value = 9250
$ws.Range("K126").Value = 85721454
$ws.Range("L126").Value = 27750
$ws.Range("M126").Value = -85718984
$ws.Range("N126").Value = -32690

$ws.Range("H141").Value = 26743.75
$ws.Range("J141").Value = 26743.75
$ws.Range("L141").Value = 26743.75
$ws.Range("N141").Value = -37103.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 4474.074
$ws.Range("J112").Value = 4625
$ws.Range("L112").Value = 13875
$ws.Range("N112").Value = -16091

$ws.Range("H121").Value = 1836.5964
$ws.Range("J121").Value = 1864.0358
$ws.Range("L121").Value = 5592.107400000001
$ws.Range("N121").Value = -8212.107400000001

$ws.Range("H122").Value = 3289.1167
$ws.Range("I122").Value = 570.6
$ws.Range("J122").Value = 3536.2546
$ws.Range("K122").Value = 5135.400000000001
$ws.Range("L122").Value = 31826.2914
$ws.Range("M122").Value = -2685.400000000001
$ws.Range("N122").Value = -36726.2914

$ws.Range("H131").Value = 6411040.5
$ws.Range("J131").Value = 811.2432
$ws.Range("L131").Value = 2433.7296
$ws.Range("N131").Value = -12513.7296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 31039.285
$ws.Range("J46").Value = 31039.285
$ws.Range("L46").Value = 31039.285
$ws.Range("N46").Value = -31351.285

$ws.Range("H102").Value = 1700.4048
$ws.Range("I102").Value = 1319
$ws.Range("J102").Value = 2551.2307
$ws.Range("K102").Value = 1319
$ws.Range("L102").Value = 2551.2307
$ws.Range("M102").Value = 303
$ws.Range("N102").Value = -5795.2307

$ws.Range("H137").Value = 42751.668
$ws.Range("J137").Value = 42751.668
$ws.Range("L137").Value = 42751.668
$ws.Range("N137").Value = -52951.668

$ws.Range("H139").Value = 59995
$ws.Range("J139").Value = 59995
$ws.Range("L139").Value = 59995
$ws.Range("N139").Value = -70275

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5573.528
$ws.Range("I40").Value = 3461.524
$ws.Range("J40").Value = 8530.333000000001
$ws.Range("K40").Value = 3461.524
$ws.Range("L40").Value = 8530.333000000001
$ws.Range("M40").Value = -3325.524
$ws.Range("N40").Value = -8802.333000000001
